# "Loop ostanich listu, komentare"
# Rename the single existing sheet to "Fiala", then add two more sheets:
#   "Bereko" - a copy of Fiala's header + first 13 employee rows
#   "Auta"   - a brand-new, empty sheet
# Finally restore each sheet's remembered cursor/selection cell.

$wb = $excel.ActiveWorkbook
$fiala = $wb.ActiveSheet
$fiala.Name = "Fiala"

# Add the two extra sheets right after "Fiala", in order, so the final tab
# order is Fiala, Bereko, Auta.
$bereko = $wb.Worksheets.Add($null, $fiala)
$bereko.Name = "Bereko"

$auta = $wb.Worksheets.Add($null, $bereko)
$auta.Name = "Auta"

# "Bereko" is a copy of the first block of "Fiala" (header row + the first
# 13 people), pasted so it starts at A1 instead of A2.
$fiala.Range("A2:J15").Copy($bereko.Range("A1"))

# Restore the selection/active-cell on every sheet.
$bereko.Range("D20").Select() | Out-Null
$auta.Range("A1").Select() | Out-Null
$fiala.Range("M21").Select() | Out-Null
